$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 329-330; this shifts the existing rows 329-400
# down to 331-402 (the former rows 399-400 fall off the end of the old
# block and become the new rows 401-402).
$ws.Rows("329:330").Insert()

# Row 329: new "Pintón" quote for the week of 2021-11-11 (serial 44511)
$ws.Range("A329").Value = 7
$ws.Range("B329").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C329").Value = "Ñuble"
$ws.Range("D329").Value = 44511
$ws.Range("E329").Value = 16
$ws.Range("F329").Value = "Fruta"
$ws.Range("G329").Value = 100108
$ws.Range("H329").Value = "Tropicales y subtropicales"
$ws.Range("I329").Value = 100108006
$ws.Range("J329").Value = "Plátano"
$ws.Range("K329").Value = "Sin especificar"
$ws.Range("L329").Value = "Pintón"
$ws.Range("M329").Value = 160
$ws.Range("N329").Value = 15500
$ws.Range("O329").Value = 16000
$ws.Range("P329").Value = 15750
$ws.Range("Q329").Value = "$/caja 20 kilos"
$ws.Range("R329").Value = "Ecuador"
$ws.Range("S329").Value = 788
$ws.Range("T329").Value = 20

# Row 330: new "Primera Pintón" quote for the same week
$ws.Range("A330").Value = 7
$ws.Range("B330").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C330").Value = "Ñuble"
$ws.Range("D330").Value = 44511
$ws.Range("E330").Value = 16
$ws.Range("F330").Value = "Fruta"
$ws.Range("G330").Value = 100108
$ws.Range("H330").Value = "Tropicales y subtropicales"
$ws.Range("I330").Value = 100108006
$ws.Range("J330").Value = "Plátano"
$ws.Range("K330").Value = "Sin especificar"
$ws.Range("L330").Value = "Primera Pintón"
$ws.Range("M330").Value = 240
$ws.Range("N330").Value = 17000
$ws.Range("O330").Value = 18000
$ws.Range("P330").Value = 17500
$ws.Range("Q330").Value = "$/caja 20 kilos"
$ws.Range("R330").Value = "Ecuador"
$ws.Range("S330").Value = 875
$ws.Range("T330").Value = 20
